$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.423999999999999
$ws.Range("C8").Value = -11.953
$ws.Range("C10").Value = -13.08
$ws.Range("D11").Value = -7.379
$ws.Range("C12").Value = -11.466
$ws.Range("D12").Value = -7.775
$ws.Range("D15").Value = -8.221
$ws.Range("D17").Value = -8.228999999999999
$ws.Range("C18").Value = -12.547
$ws.Range("C25").Value = -12.141
$ws.Range("D26").Value = -7.255000000000001
$ws.Range("D27").Value = -7.572
$ws.Range("D28").Value = -7.609
$ws.Range("D32").Value = -7.362
$ws.Range("C37").Value = -13.608
$ws.Range("D37").Value = -7.361
$ws.Range("D41").Value = -7.578999999999999
$ws.Range("D47").Value = -7.452
$ws.Range("D51").Value = -8.022
$ws.Range("C55").Value = -13.837
$ws.Range("D65").Value = -7.741
$ws.Range("C68").Value = -11.134
$ws.Range("D73").Value = -7.374000000000001
$ws.Range("C77").Value = -13.557
$ws.Range("C78").Value = -13.401
$ws.Range("C79").Value = -12.603
$ws.Range("C80").Value = -12.858
$ws.Range("C81").Value = -13.333
$ws.Range("C82").Value = -12.354
$ws.Range("C84").Value = -12.9
$ws.Range("D84").Value = -7.931
$ws.Range("D85").Value = -8.654
$ws.Range("D89").Value = -8.327
$ws.Range("D93").Value = -7.024000000000001
$ws.Range("D95").Value = -7.506
$ws.Range("D98").Value = -7.231
$ws.Range("D99").Value = -8.215
$ws.Range("C101").Value = -12.377
$ws.Range("D101").Value = -7.723999999999999
$ws.Range("C102").Value = -13.327
$ws.Range("D102").Value = -7.279999999999999
